$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cyclic rotation of data rows 2, 3, 4:
#   new row2 <= old row4
#   new row3 <= old row2
#   new row4 <= old row3
# (read with Value2 first since these rows overlap in the rotation)

$oldA2 = $ws.Range("A2").Value2
$oldB2 = $ws.Range("B2").Value2
$oldD2 = $ws.Range("D2").Value2
$oldE2 = $ws.Range("E2").Value2
$oldF2 = $ws.Range("F2").Value2
$oldG2 = $ws.Range("G2").Value2
$oldH2 = $ws.Range("H2").Value2
$oldQ2 = $ws.Range("Q2").Value2
$oldR2 = $ws.Range("R2").Value2
$oldAC2 = $ws.Range("AC2").Value2

$oldA3 = $ws.Range("A3").Value2
$oldB3 = $ws.Range("B3").Value2
$oldD3 = $ws.Range("D3").Value2
$oldE3 = $ws.Range("E3").Value2
$oldF3 = $ws.Range("F3").Value2
$oldG3 = $ws.Range("G3").Value2
$oldH3 = $ws.Range("H3").Value2
$oldQ3 = $ws.Range("Q3").Value2
$oldR3 = $ws.Range("R3").Value2
$oldAC3 = $ws.Range("AC3").Value2

$oldA4 = $ws.Range("A4").Value2
$oldB4 = $ws.Range("B4").Value2
$oldD4 = $ws.Range("D4").Value2
$oldE4 = $ws.Range("E4").Value2
$oldF4 = $ws.Range("F4").Value2
$oldG4 = $ws.Range("G4").Value2
$oldH4 = $ws.Range("H4").Value2
$oldQ4 = $ws.Range("Q4").Value2
$oldR4 = $ws.Range("R4").Value2
$oldAC4 = $ws.Range("AC4").Value2

# Row 2 becomes old row 4
$ws.Range("A2").Value = $oldA4
$ws.Range("B2").Value = $oldB4
$ws.Range("D2").Value = $oldD4
$ws.Range("E2").Value = $oldE4
$ws.Range("F2").Value = $oldF4
$ws.Range("G2").Value = $oldG4
$ws.Range("H2").Value = $oldH4
$ws.Range("Q2").Value = $oldQ4
$ws.Range("R2").Value = $oldR4
$ws.Range("AC2").Value = $oldAC4

# Row 3 becomes old row 2
$ws.Range("A3").Value = $oldA2
$ws.Range("B3").Value = $oldB2
$ws.Range("D3").Value = $oldD2
$ws.Range("E3").Value = $oldE2
$ws.Range("F3").Value = $oldF2
$ws.Range("G3").Value = $oldG2
$ws.Range("H3").Value = $oldH2
$ws.Range("Q3").Value = $oldQ2
$ws.Range("R3").Value = $oldR2
$ws.Range("AC3").Value = $oldAC2

# Row 4 becomes old row 3
$ws.Range("A4").Value = $oldA3
$ws.Range("B4").Value = $oldB3
$ws.Range("D4").Value = $oldD3
$ws.Range("E4").Value = $oldE3
$ws.Range("F4").Value = $oldF3
$ws.Range("G4").Value = $oldG3
$ws.Range("H4").Value = $oldH3
$ws.Range("Q4").Value = $oldQ3
$ws.Range("R4").Value = $oldR3
$ws.Range("AC4").Value = $oldAC3
